$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

$data = @(
    ,@('Test Case ID', 'Test scenario', 'Status')
    ,@('iEX_TC_ID_10', '@iExamRegression Validation of Sign out at Exam Start Page', 'passed')
    ,@('iEX_TC_ID_9', '@iExamRegression Validation of user authentication by valid Candidate Creadentials', 'passed')
    ,@('iEX_TC_ID_16', '@iExamRegression Validation of Sign out at Candidate Dashboard Page', 'passed')
    ,@('iEX_TC_ID_14', '@iExamRegression Validate login with inactive user', 'passed')
    ,@('iEX_TC_ID_4', '@iExamRegression Validation of Empty Field validation', 'passed')
    ,@('iEX_TC_ID_19', '@iExamRegression Validation of Candidate App Dashboard Details', 'passed')
    ,@('iEX_TC_ID_20', '@iExamRegression Validation of Candidate Start Exam', 'passed')
    ,@('iEX_TC_ID_22,iEX_TC_ID_78', '@iExamRegression Password will not be asked after clicking on "Start Exam"', 'passed')
    ,@('iEX_TC_ID_18', '@iExamRegression Validation of multiple candidate trying to login to same Exam', 'passed')
    ,@('iEX_TC_ID_31', '@iExamRegression Verify Validation of Browser back button on Candidate Dashboard', 'passed')
    ,@('iEX_TC_ID_34', '@iExamRegression Verify Validation of  function keys on Candidate Dashboard', 'passed')
    ,@('iEX_TC_ID_36', '@iExamRegression Verify Validation of Changing Font Size to increase on the Dashboard', 'passed')
    ,@('iEX_TC_ID_33', '@iExamRegression Verify Validation of Browser Reload option on Candidate Dashboard', 'passed')
    ,@('iEX_TC_ID_36B', '@iExamRegression Verify Validation of Changing Font Size to Decrease on the Dashboard', 'passed')
    ,@('iEX_TC_ID_45', '@iExamRegression Verify Validation of the "Timer" - After clicking on "Start Exam" Button', 'passed')
    ,@('iEX_TC_ID_32,iEX_TC_ID_42', '@iExamRegression Verify Validation of Browser back button after Exam Starts', 'passed')
    ,@('iEX_TC_ID_41', '@iExamRegression Verify Validation of Browser Reload option after Exam Starts', 'passed')
    ,@('iEX_TC_ID_43', '@iExamRegression Verify Validation of  function keys after Exam Starts', 'passed')
    ,@('iEX_TC_ID_225', '@iExamRegression Verify Candidate while attending exam - Offline - try to press any F1 to F12 keys', 'passed')
    ,@('iEX_TC_ID_56,iEX_TC_ID_212', '@iExamRegression Verify Validation of abrupt closure of browser', 'passed')
    ,@('iEX_TC_ID_50', '@iExamRegression Verify Validation of Updated cloud symbol', 'passed')
    ,@('iEX_TC_ID_74', '@iExamRegression Verify Validation of Flag for Review', 'passed')
    ,@('iEX_TC_ID_55,iEX_TC_ID_238', '@iExamRegression Verify Validation of the following in the Exam Section
Exam name, Candidate name, Exam timer, no of questions as per exam created EluminaRegressioniExamTC_051', 'passed')
    ,@('iEX_TC_ID_59', '@iExamRegression Verify candidate adds and saves Notes EluminaRegressioniExamTC_038', 'passed')
    ,@('iEX_TC_ID_56B', '@iExamRegression Verify Candidate is able to close the broswer and answered questions are saved after logging in again', 'passed')
    ,@('iEX_TC_ID_83', '@iExamRegression Validation of candidate able to submit the exam when not answering all questions', 'passed')
    ,@('iEX_TC_ID_88', '@iExamRegression Validation of "Time Remaining" pop-up when the just before the exam time runs out', 'passed')
    ,@('iEX_TC_ID_62', '@iExamRegression Verify Validation of using Calculator EluminaRegressioniExamTC_039 ', 'passed')
    ,@('iEX_TC_ID_97,iEX_TC_ID_98,iEX_TC_ID_99,iEX_TC_ID_101,iEX_TC_ID_102,iEX_TC_ID_103,iEX_TC_ID_105', '@iExamRegression Validation of each component for all type of questions in Candidate flow.', 'failed')
    ,@('iEX_TC_ID_64', '@iExamRegression Verify Validation of using Highlighter EluminaRegressioniExamTC_040 ', 'passed')
    ,@('iEX_TC_ID_68', '@iExamRegression Verify Validation of Candidate attends All Question type', 'passed')
    ,@('iEX_TC_ID_99', '@iExamRegression Verify Validation of VSAQ Question type EluminaRegressioniExamTC_043', 'passed')
    ,@('iEX_TC_ID_49', '@iExamRegression Validation of Exam content page ->Timer verification in all pages ', 'passed')
    ,@('iEX_TC_ID_54', '@iExamRegression Validation of Exam section page.', 'passed')
    ,@('iEX_TC_ID_61', '@iExamRegression Validation of Exam Section > Font Size ', 'passed')
    ,@('iEX_TC_ID_48', '@iExamRegression Verify Timer Expires Alert Box  in Exam Content page', 'passed')
    ,@('iEX_TC_ID_69,iEX_TC_ID_240', '@iExamRegression Validation of Exam Section > Question and Answer save navigation validation using Previous / Next buttons.', 'passed')
    ,@('iEX_TC_ID_70', '@iExamRegression Validate the Exam sheet where the Question numbers are displayed are Gray by Default', 'passed')
    ,@('iEX_TC_ID_64', '@iExamRegression Validation of Exam Section > Highlighter tool highlights save scenario 1', 'passed')
    ,@('iEX_TC_ID_71', '@iExamRegression Verify Validate the Exam sheet where the Question numbers are displayed in Red, when Flag for Review is selected
Verification of Timer Expires Alert Box  in Exam Content page', 'passed')
    ,@('iEX_TC_ID_30', '@iExamRegression Verify Validation of Candidate dashboard - Exam scheduled Date / Time is over EluminaRegressioniExamTC_024', 'timedOut')
    ,@('iEX_TC_ID_79', '@iExamRegression Verify Validation of Candidate attends All Question type in offline', 'passed')
    ,@('iEX_TC_ID_76', '@iExamRegression Validation of Exam section > Candidate attend the exam in normal flow ', 'passed')
    ,@('iEX_TC_ID_30B', '@iExamRegression Verify Validation of Candidate dashboard - Exam scheduled Date / Time is over EluminaRegressioniExamTC_024A', 'passed')
    ,@('iEX_TC_ID_72', '@iExamRegression Validate the Exam sheet where the Question numbers are displayed in Orange when In Progress', 'passed')
    ,@('iEX_TC_ID_39', '@iExamRegression Verify Validation of Exam content section page contents and details', 'passed')
    ,@('iEX_TC_ID_73', '@iExamRegression Verify Validate the Exam sheet where the Question numbers are displayed in Pink, when Notes are added while answering ', 'passed')
    ,@('iEX_TC_ID_45', '@iExamRegression Verify Validation of Exam content page - Timer Validation (Before Timer Expires)', 'passed')
    ,@('iEX_TC_ID_47', '@iExamRegression Verify Validation of Exam content page  - Check box', 'passed')
    ,@('iEX_TC_ID_46', '@iExamRegression Verify Validation of Exam content page - Scroll bar', 'passed')
    ,@('iEX_TC_ID_127', '@iExamRegression Verify Validation of login for invigilator for Invalid credential', 'passed')
    ,@('iEX_TC_ID_126', '@iExamRegression Verify Validation of login for invigilator for valid credential', 'passed')
    ,@('iEX_TC_ID_95', '@iExamRegression Verify Validation of Candidate attends All Question type in offline', 'passed')
    ,@('iEX_TC_ID_152', '@iExamRegression Verify Invigilator marks attendance for All candidates', 'passed')
    ,@('iEX_TC_ID_148', '@iExamRegression Validation of Invigilator Dashboard after the Exam Completion by candidate (once the time exceeds)', 'passed')
    ,@('iEX_TC_ID_25D', '@iExamRegression Verify Login Application for PM', 'passed')
    ,@('iEX_TC_ID_148B', '@iExamRegression Verify Validation of Invigilator Dashboard ', 'passed')
    ,@('iEX_TC_ID_199', '@iExamRegression Validate auto response when chat is initiated out of working hours.', 'passed')
    ,@('iEX_TC_ID_235', '@iExamRegression Validation of Candidate dashboard > with Unassigned / Cancelled Exam ', 'passed')
    ,@('iEX_TC_ID_242,243,iEX_TC_ID_244,iEX_TC_ID_234,iEX_TC_ID_241', '@iExamRegression Verify Validation of Candidate attend Practice Exam', 'timedOut')
    ,@('iEX_TC_ID_236', '@iExamRegression Validation of Candidate dashboard > Exam scheduled Date / Time is over.', 'failed')
    ,@('iEX_TC_ID_90,iEX_TC_ID_96', '@iExamRegression Verify Validation of Survey screen and Validate Survey screen where candidate
    can provide feedback in comment section TC-084 and TC-085', 'passed')
    ,@('iEX_TC_ID_91', '@iExamRegression Survey section > validate exam name, candidate name, exam timer, no of questions is as per exam settings.', 'passed')
    ,@('iEX_TC_ID_92', '@iExamRegression Survey section > Flag for review.', 'passed')
    ,@('iEX_TC_ID_57', '@iExamRegression Validation of Exam Section> Image viewer', 'timedOut')
    ,@('iEX_TC_ID_032', '@iExamRegression Validation of Sign out at Exam Start Page', 'passed')
    ,@('iEX_TC_ID_77', '@iExamRegression Validation of Exam section > Candidate attend the exam by answering the questions randomly (Example: Monkey Jump)', 'passed')
    ,@('iEX_TC_ID_93', '@iExamRegression Validation of Survey section > Before Exam Time Expires', 'passed')
    ,@('iEX_TC_ID_89', '@iExamRegression Validation of Exam section > After Exam Time Expires', 'passed')
    ,@('iEX_TC_ID_3', '@LowPriorityiExamCases Validation of Client Logo and name, color, font size and font type of different elements in the login page', 'passed')
    ,@('iEX_TC_ID_5', '@LowPriorityiExamCases Verify Validation of User Id Field', 'passed')
    ,@('iEX_TC_ID_11', '@LowPriorityiExamCases Validation of Login using Invalid credentials Scenario 1 TC-09', 'passed')
    ,@('iEX_TC_ID_6', '@LowPriorityiExamCases Verify Validation of Password popup Field', 'passed')
    ,@('iEX_TC_ID_12,iEX_TC_ID_13', '@LowPriorityiExamCases Validation of Login using Invalid credentials Scenario 2 TC-10', 'passed')
    ,@('iEX_TC_ID_229', '@LowPriorityiExamCases Verify Validation of Login using valid credentials in browser offline', 'passed')
    ,@('iEX_TC_ID_89A', '@iExamRegression Validation of Exam section > After Exam Time Expires ', 'failed')
    ,@('iEX_TC_ID_17', '@LowPriorityiExamCases Validation of Sign out at Candidate Dashboard Page TC-14', 'passed')
    ,@('iEX_TC_ID_192', '@LowPriorityiExamCases Validation of Candidate Dashboard page - Chat App', 'passed')
    ,@('iEX_TC_ID_17', '@LowPriorityiExamCases Validation of Sign out at Exam Start Page TC-15', 'passed')
    ,@('iEX_TC_ID_191', '@LowPriorityiExamCases Validation of Candidate Login page', 'passed')
    ,@('iEX_TC_ID_60', '@LowPriorityiExamCases Verify candidate adds and saves Notes EluminaRegressioniExamTC_038', 'passed')
    ,@('iEX_TC_ID_193', '@LowPriorityiExamCases Validation of Exam content section >  Chat App', 'passed')
    ,@('iEX_TC_ID_36', '@LowPriorityiExamCases Validation of Candidate dashboard > Exam Start page Font size validation', 'passed')
    ,@('iEX_TC_ID_94', '@iExamRegression Validation of Survey section > After Exam Time Expires', 'passed')
    ,@('iEX_TC_ID_52', '@LowPriorityiExamCases Validation of Exam content page -> Questions download at content section-> Browser shut down', 'passed')
    ,@('iEX_TC_ID_230', '@LowPriorityiExamCases Validation of copy paste in password field', 'passed')
    ,@('iEX_TC_ID_194', '@LowPriorityiExamCases Verify Validation of Exam section page  > Chat App', 'passed')
    ,@('iEX_TC_ID_44', '@LowPriorityiExamCases Validation of Exam content page -> (Terms & Conditions) Font size validation TC-33', 'passed')
    ,@('iEX_TC_ID_52', '@LowPriorityiExamCases Validation of Exam content page -> Questions download at content section-> Browser shut down TC-40', 'passed')
    ,@('iEX_TC_ID_198', '@LowPriorityiExamCases Validate when user go offline while using chat app.', 'passed')
    ,@('iEX_TC_ID_224', '@LowPriorityiExamCases Candidate while attending exam - Online - try to press any F1 to F12 keys TC-201', 'passed')
    ,@('iEX_TC_ID_220', '@LowPriorityiExamCases Validation of Candidate Recovery Sync Page Validation', 'passed')
    ,@('iEX_TC_ID_197', '@LowPriorityiExamCases Validation of textbox capabilities of chat app', 'passed')
    ,@('iEX_TC_ID_53', '@LowPriorityiExamCases Validation of Exam content page -> Questions download at content section-> Candidate goes offline', 'passed')
    ,@('iEX_TC_ID_222', '@LowPriorityiExamCases "Validation of 
exam response recovery page with no recovery file available in Download offline recovery file"', 'failed')
    ,@('iEX_TC_ID_228', '@LowPriorityiExamCases No other key combinations could cause the candidates to logout and/or lose data TC-205', 'passed')
    ,@('iEX_TC_ID_108,iEX_TC_ID_115', '@LowPriorityiExamCases Validation of Candidate attending Exam in Online - Offline and submit in Offline Mode', 'passed')
    ,@('iEX_TC_ID_114', '@LowPriorityiExamCases Verify Validation of Candidate attends All Question type', 'passed')
    ,@('iEX_TC_ID_116', '@LowPriorityiExamCases Validation of Candidate Attending Exam in Online (Abort and Resume Exam)', 'passed')
    ,@('iEX_TC_ID_116B', '@LowPriorityiExamCases Validation Candidate Login again', 'passed')
    ,@('iEX_TC_ID_195', '@LowPriorityiExamCases Verify Validation of Review Exam page  > Chat App TC-173', 'failed')
    ,@('iEX_TC_ID_196', '@LowPriorityiExamCases Verify Validation of Submit Exam page  > Chat App TC-173', 'timedOut')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}
